$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the daily log for Ciruela
# (Feria Lagunitas de Puerto Montt). Insert a new row above row 64,
# shifting the existing rows 64-107 down to 65-108, and populate the
# new row with the latest reading.
$ws.Rows.Item(64).Insert()

$ws.Range("A64").Value = 4
$ws.Range("B64").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C64").Value = "Los Lagos"
$ws.Range("D64").Value = 44566
$ws.Range("E64").Value = 10
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100103
$ws.Range("H64").Value = "Frutos de hueso (carozo)"
$ws.Range("I64").Value = 100103002
$ws.Range("J64").Value = "Ciruela"
$ws.Range("K64").Value = "Black Amber"
$ws.Range("L64").Value = "Primera"
$ws.Range("M64").Value = 200
$ws.Range("N64").Value = 18000
$ws.Range("O64").Value = 18500
$ws.Range("P64").Value = 18250
$ws.Range("Q64").Value = "`$/caja 15 kilos granel"
$ws.Range("R64").Value = "Región de O'Higgins"
$ws.Range("S64").Value = 1217
$ws.Range("T64").Value = 15
